# ExcelSubTable interface: translate sheet contents from Russian to English
# and update selection to match new header position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "First Cell"
$ws.Range("C4").Value = "Foo"
$ws.Range("F7").Value = "Data1"
$ws.Range("F8").Value = "Data2"
$ws.Range("F9").Value = "Data3"
$ws.Range("F10").Value = "Data4"
$ws.Range("G6").Value = "Header1"
$ws.Range("H6").Value = "Header2"
$ws.Range("I6").Value = "Header3"
$ws.Range("J6").Value = "Header4"
$ws.Range("I7").Value = "Foo"
$ws.Range("I9").Value = "Bas"

$ws.Range("G6").Select()
